$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.806.76"
$ws.Range("E2").Value = "  -0.40%  "

$ws.Range("D3").Value = "2.031.51"
$ws.Range("E3").Value = "  -0.95%  "

$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.40"
$ws.Range("E5").Value = "  -1.08%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.609"
$ws.Range("E6").Value = "  -1.14%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.27"
$ws.Range("E7").Value = "  +2.15%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("E9").Value = "  -0.13%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0812"
$ws.Range("E10").Value = "  +0.59%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.103"
$ws.Range("E11").Value = "  -0.14%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.64"
$ws.Range("E12").Value = "  +0.69%  "

$ws.Range("D13").Value = "2.332.61"
$ws.Range("E13").Value = "  -1.13%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.13"
$ws.Range("E14").Value = "  +2.51%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.763"
$ws.Range("E15").Value = "  +1.45%  "

$ws.Range("E16").Value = "  -1.57%  "

$ws.Range("D17").Value = "2.034.15"
$ws.Range("E17").Value = "  -0.58%  "

$ws.Range("D18").Value = "37.785.52"
$ws.Range("E18").Value = "  -0.29%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.04"
$ws.Range("E19").Value = "  -1.50%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.98"
$ws.Range("E20").Value = "  +0.42%  "

$ws.Range("E21").Value = "  -0.71%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "223.80"
$ws.Range("E22").Value = "  -0.32%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.02%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.40"
$ws.Range("E24").Value = "  -2.27%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.21"
$ws.Range("E25").Value = "  -1.65%  "

$ws.Range("E26").Value = "  +0.09%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.07"
$ws.Range("E27").Value = "  -0.77%  "

$ws.Range("E28").Value = "  -2.34%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.96"
$ws.Range("E29").Value = "  -0.20%  "

$ws.Range("E30").Value = "  -4.34%  "

$ws.Range("E31").Value = "  +1.05%  "

$ws.Range("E32").Value = "  -1.80%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.08"
$ws.Range("E33").Value = "  +2.71%  "

$ws.Range("E34").Value = "  -1.50%  "

$ws.Range("E35").Value = "  -1.35%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.38"
$ws.Range("E36").Value = "  +6.97%  "

$ws.Range("E37").Value = "  -3.46%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.25"
$ws.Range("E38").Value = "  -1.78%  "

$ws.Range("E39").Value = "  -0.22%  "

$ws.Range("D40").Value = "1.526.30"
$ws.Range("E40").Value = "  +2.91%  "

$ws.Range("E41").Value = "  +0.54%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "96.76"
$ws.Range("E42").Value = "  -0.91%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.80"
$ws.Range("E43").Value = "  +0.77%  "

$ws.Range("E44").Value = "  -0.61%  "

$ws.Range("E45").Value = "  -1.61%  "

$ws.Range("E46").Value = "  -1.43%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.04"
$ws.Range("E47").Value = "  -4.29%  "

$ws.Range("E48").Value = "  -0.72%  "

$ws.Range("E49").Value = "  -0.27%  "

$ws.Range("E50").Value = "  +0.98%  "

$ws.Range("D51").Value = "2.220.79"
$ws.Range("E51").Value = "  -1.16%  "
